$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5655645728111267
$ws.Range("B1").Value = 2.045938014984131
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.275944590568542
$ws.Range("E1").Value = 1.313208937644958
